$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.320746
$ws.Range("H2").Value = 288.962238
$ws.Range("I2").Value = 0.3809824610908788
$ws.Range("J2").Value = 0.3809824610908788
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 15.35884066666667
$ws.Range("N2").Value = 46.076522
$ws.Range("O2").Value = 0.1012042817263867
$ws.Range("P2").Value = 0.1012042817263867
$ws.Range("Q2").Value = 1479.37499070847
$ws.Range("R2").Value = 13314.37491637623
$ws.Range("S2").Value = 0.03855705632505344
$ws.Range("T2").Value = 0.03855705632505344

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.320746
$ws.Range("H3").Value = 288.962238
$ws.Range("I3").Value = 0.3809824610908788
$ws.Range("J3").Value = 0.3809824610908788
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 50.59256466666667
$ws.Range("N3").Value = 151.777694
$ws.Range("O3").Value = 0.3333704853712116
$ws.Range("P3").Value = 0.3333704853712116
$ws.Range("Q3").Value = 4873.113570746575
$ws.Range("R3").Value = 43858.02213671917
$ws.Range("S3").Value = 0.127008307971785
$ws.Range("T3").Value = 0.127008307971785

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.320746
$ws.Range("H4").Value = 288.962238
$ws.Range("I4").Value = 0.3809824610908788
$ws.Range("J4").Value = 0.3809824610908788
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 60.37715666666667
$ws.Range("N4").Value = 181.13147
$ws.Range("O4").Value = 0.397844271305776
$ws.Range("P4").Value = 0.397844271305776
$ws.Range("Q4").Value = 5815.572771492207
$ws.Range("R4").Value = 52340.15494342986
$ws.Range("S4").Value = 0.1515716896129818
$ws.Range("T4").Value = 0.1515716896129818

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.320746
$ws.Range("H5").Value = 288.962238
$ws.Range("I5").Value = 0.3809824610908788
$ws.Range("J5").Value = 0.3809824610908788
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 25.43221733333333
$ws.Range("N5").Value = 76.29665199999999
$ws.Range("O5").Value = 0.1675809615966257
$ws.Range("P5").Value = 0.1675809615966258
$ws.Range("Q5").Value = 2449.650145980797
$ws.Range("R5").Value = 22046.85131382717
$ws.Range("S5").Value = 0.06384540718105852
$ws.Range("T5").Value = 0.06384540718105854

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.46467533333333
$ws.Range("H6").Value = 55.394026
$ws.Range("I6").Value = 0.07303429161291354
$ws.Range("J6").Value = 0.07303429161291354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.35884066666667
$ws.Range("N6").Value = 46.076522
$ws.Range("O6").Value = 0.1012042817263867
$ws.Range("P6").Value = 0.1012042817263867
$ws.Range("Q6").Value = 283.5960064063968
$ws.Range("R6").Value = 2552.364057657572
$ws.Range("S6").Value = 0.00739138302408038
$ws.Range("T6").Value = 0.007391383024080381

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.46467533333333
$ws.Range("H7").Value = 55.394026
$ws.Range("I7").Value = 0.07303429161291354
$ws.Range("J7").Value = 0.07303429161291354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 50.59256466666667
$ws.Range("N7").Value = 151.777694
$ws.Range("O7").Value = 0.3333704853712116
$ws.Range("P7").Value = 0.3333704853712116
$ws.Range("Q7").Value = 934.1752808506715
$ws.Range("R7").Value = 8407.577527656043
$ws.Range("S7").Value = 0.02434747724373959
$ws.Range("T7").Value = 0.02434747724373959

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.46467533333333
$ws.Range("H8").Value = 55.394026
$ws.Range("I8").Value = 0.07303429161291354
$ws.Range("J8").Value = 0.07303429161291354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 60.37715666666667
$ws.Range("N8").Value = 181.13147
$ws.Range("O8").Value = 0.397844271305776
$ws.Range("P8").Value = 0.397844271305776
$ws.Range("Q8").Value = 1114.844595399802
$ws.Range("R8").Value = 10033.60135859822
$ws.Range("S8").Value = 0.02905627452707314
$ws.Range("T8").Value = 0.02905627452707314

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.46467533333333
$ws.Range("H9").Value = 55.394026
$ws.Range("I9").Value = 0.07303429161291354
$ws.Range("J9").Value = 0.07303429161291354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.43221733333333
$ws.Range("N9").Value = 76.29665199999999
$ws.Range("O9").Value = 0.1675809615966257
$ws.Range("P9").Value = 0.1675809615966258
$ws.Range("Q9").Value = 469.5976360667723
$ws.Range("R9").Value = 4226.378724600952
$ws.Range("S9").Value = 0.01223915681802043
$ws.Range("T9").Value = 0.01223915681802043

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.909391
$ws.Range("H10").Value = 368.728173
$ws.Range("I10").Value = 0.4861499128584522
$ws.Range("J10").Value = 0.4861499128584522
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.35884066666667
$ws.Range("N10").Value = 46.076522
$ws.Range("O10").Value = 0.1012042817263867
$ws.Range("P10").Value = 0.1012042817263867
$ws.Range("Q10").Value = 1887.745752806034
$ws.Range("R10").Value = 16989.7117752543
$ws.Range("S10").Value = 0.04920045274218511
$ws.Range("T10").Value = 0.04920045274218512

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.909391
$ws.Range("H11").Value = 368.728173
$ws.Range("I11").Value = 0.4861499128584522
$ws.Range("J11").Value = 0.4861499128584522
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 50.59256466666667
$ws.Range("N11").Value = 151.777694
$ws.Range("O11").Value = 0.3333704853712116
$ws.Range("P11").Value = 0.3333704853712116
$ws.Range("Q11").Value = 6218.301312308117
$ws.Range("R11").Value = 55964.71181077306
$ws.Range("S11").Value = 0.1620680324127944
$ws.Range("T11").Value = 0.1620680324127944

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 122.909391
$ws.Range("H12").Value = 368.728173
$ws.Range("I12").Value = 0.4861499128584522
$ws.Range("J12").Value = 0.4861499128584522
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 60.37715666666667
$ws.Range("N12").Value = 181.13147
$ws.Range("O12").Value = 0.397844271305776
$ws.Range("P12").Value = 0.397844271305776
$ws.Range("Q12").Value = 7420.91955621159
$ws.Range("R12").Value = 66788.27600590431
$ws.Range("S12").Value = 0.1934119578265374
$ws.Range("T12").Value = 0.1934119578265374

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 122.909391
$ws.Range("H13").Value = 368.728173
$ws.Range("I13").Value = 0.4861499128584522
$ws.Range("J13").Value = 0.4861499128584522
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 25.43221733333333
$ws.Range("N13").Value = 76.29665199999999
$ws.Range("O13").Value = 0.1675809615966257
$ws.Range("P13").Value = 0.1675809615966258
$ws.Range("Q13").Value = 3125.858344219643
$ws.Range("R13").Value = 28132.72509797679
$ws.Range("S13").Value = 0.08146946987693522
$ws.Range("T13").Value = 0.08146946987693524

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.127183
$ws.Range("H14").Value = 45.381549
$ws.Range("I14").Value = 0.05983333443775553
$ws.Range("J14").Value = 0.05983333443775553
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 15.35884066666667
$ws.Range("N14").Value = 46.076522
$ws.Range("O14").Value = 0.1012042817263867
$ws.Range("P14").Value = 0.1012042817263867
$ws.Range("Q14").Value = 232.3359934325087
$ws.Range("R14").Value = 2091.023940892578
$ws.Range("S14").Value = 0.006055389635067723
$ws.Range("T14").Value = 0.006055389635067724

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.127183
$ws.Range("H15").Value = 45.381549
$ws.Range("I15").Value = 0.05983333443775553
$ws.Range("J15").Value = 0.05983333443775553
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 50.59256466666667
$ws.Range("N15").Value = 151.777694
$ws.Range("O15").Value = 0.3333704853712116
$ws.Range("P15").Value = 0.3333704853712116
$ws.Range("Q15").Value = 765.3229841520007
$ws.Range("R15").Value = 6887.906857368006
$ws.Range("S15").Value = 0.01994666774289259
$ws.Range("T15").Value = 0.01994666774289259

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.127183
$ws.Range("H16").Value = 45.381549
$ws.Range("I16").Value = 0.05983333443775553
$ws.Range("J16").Value = 0.05983333443775553
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 60.37715666666667
$ws.Range("N16").Value = 181.13147
$ws.Range("O16").Value = 0.397844271305776
$ws.Range("P16").Value = 0.397844271305776
$ws.Range("Q16").Value = 913.3362979163368
$ws.Range("R16").Value = 8220.026681247031
$ws.Range("S16").Value = 0.02380434933918364
$ws.Range("T16").Value = 0.02380434933918364

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.127183
$ws.Range("H17").Value = 45.381549
$ws.Range("I17").Value = 0.05983333443775553
$ws.Range("J17").Value = 0.05983333443775553
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 25.43221733333333
$ws.Range("N17").Value = 76.29665199999999
$ws.Range("O17").Value = 0.1675809615966257
$ws.Range("P17").Value = 0.1675809615966258
$ws.Range("Q17").Value = 384.7178056971053
$ws.Range("R17").Value = 3462.460251273948
$ws.Range("S17").Value = 0.01002692772061157
$ws.Range("T17").Value = 0.01002692772061158
